# Apply the weekly report update for WR_89787325_WeekEnding_072725.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: report generation timestamp ---
$ws.Range("D5").Value = "Report Generated On: 08/26/2025 09:59 AM"

# --- Summary: total billed amount ---
$ws.Range("C8").Value = 13242.56

# --- Scope ID # is cleared out ---
$ws.Range("G10").Value = ""

# --- Thursday (07/24/2025) detail block, rows 16-20, plus TOTAL row 21 ---
$ws.Range("H16").Value = 478.55
$ws.Range("H17").Value = 238.2
$ws.Range("H18").Value = 476.4
$ws.Range("H19").Value = 1429.2
$ws.Range("H20").Value = 1429.2
$ws.Range("H21").Value = 4051.55

# --- Friday (07/25/2025) detail block, rows 26-49, plus TOTAL row 50 ---
$ws.Range("H26").Value = 31.72
$ws.Range("H27").Value = 63.44
$ws.Range("H28").Value = 188.34
$ws.Range("H29").Value = 94.17
$ws.Range("H30").Value = 55.18
$ws.Range("H31").Value = 62.42
$ws.Range("H32").Value = 31.72
$ws.Range("H33").Value = 63.44
$ws.Range("H34").Value = 282.51
$ws.Range("H35").Value = 94.17
$ws.Range("H36").Value = 55.18
$ws.Range("H37").Value = 62.42
$ws.Range("H38").Value = 94.17
$ws.Range("H39").Value = 55.18
$ws.Range("H40").Value = 62.42
$ws.Range("H41").Value = 31.72
$ws.Range("H42").Value = 63.44
$ws.Range("H43").Value = 282.51
$ws.Range("H44").Value = 94.17
$ws.Range("H45").Value = 62.42
$ws.Range("H46").Value = 648.53
$ws.Range("H47").Value = 648.53
$ws.Range("H48").Value = 648.53
$ws.Range("H49").Value = 648.53
$ws.Range("H50").Value = 4424.86

# --- Saturday (07/26/2025) detail block, rows 55-57, plus TOTAL row 58 ---
$ws.Range("H55").Value = 478.55
$ws.Range("H56").Value = 476.4
$ws.Range("H57").Value = 3811.2
$ws.Range("H58").Value = 4766.15
